# Update column G ("K") values on the active sheet.
# New values per row, computed/regenerated (per commit message: "regen
# save_data to use K instead of Strike#, regen std/mean, calc and write
# s_vals").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 2;  3 = 2;  4 = 1;  5 = 1;  6 = 0;  7 = 1;  8 = 0;  9 = 1;  10 = 1;
    11 = 2; 12 = 2; 13 = 0; 14 = 2; 15 = 2; 16 = 2; 17 = 2; 18 = 0; 19 = 0;
    20 = 1; 21 = 1; 22 = 0; 23 = 2; 24 = 0; 25 = 3; 26 = 0; 27 = 3; 28 = 0;
    29 = 0; 30 = 0; 31 = 0; 32 = 0; 33 = 0; 34 = 1; 35 = 0; 36 = 2; 37 = 0;
    38 = 0; 39 = 1; 40 = 1; 41 = 2; 42 = 1; 43 = 1; 44 = 1; 45 = 1; 46 = 0;
    47 = 1; 48 = 1; 49 = 2; 50 = 1; 51 = 0; 52 = 1; 53 = 1; 54 = 2; 55 = 1;
    56 = 1; 57 = 1; 58 = 1; 59 = 3; 60 = 1; 61 = 0; 62 = 1; 63 = 2; 64 = 2;
    65 = 2; 66 = 1; 67 = 0; 68 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
